$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the flight time (I7) from 29 to 32 (dependent formulas recalc automatically)
$ws.Range("I7").Value = 32

# Update active cell selection to I8
$ws.Range("I8").Select()
